# Photo viewer now displays resolution, date added, and date modified.
#
# R-PV3 "Display timestamp of added photos and modified" and R-PV4 "Display
# resolution" move from In Progress / Incomplete to Complete, each with a
# completion date of 11/3/2011 (serial 40850) in column C - matching the
# other completed requirement rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 32 : R-PV3 ("Display timestamp of added photos and modified") ---
# Status: "In Progress" (yellow) -> "Complete" (green)
$ws.Range("B32").Value = "Complete"
$ws.Range("B32").Font.Color = 5287936   # RGB(0,176,80) == style used by other "Complete" cells

# Date column: copy the date-format from an existing date cell, then set the value
$ws.Range("C2").Copy()
$ws.Range("C32").PasteSpecial(-4122)    # xlPasteFormats
$ws.Range("C32").Value = 40850

# --- Row 33 : R-PV4 ("Display resolution") ---
# Status: "Incomplete" (red) -> "Complete" (green)
$ws.Range("B33").Value = "Complete"
$ws.Range("B33").Font.Color = 5287936

$ws.Range("C2").Copy()
$ws.Range("C33").PasteSpecial(-4122)    # xlPasteFormats
$ws.Range("C33").Value = 40850

# Restore the active selection to where the edit was made
$ws.Range("C29").Select()
